$d = $word.ActiveDocument

# --- Edit 1: paragraph "...brush off........know-how" -> "... strenghten........know-how"
$p1 = $d.Paragraphs.Item(3)
$r1 = $p1.Range
$r1.Find.ClearFormatting()
$r1.Find.Execute("brush off", $true, $false, $false, $false, $false, $true, 1, $false, " strenghten", 2)

# --- Edit 2: paragraph "...enhance........interpersonal skills" -> "...broaden........interpersonal skills"
$p2 = $d.Paragraphs.Item(5)
$r2 = $p2.Range
$r2.Find.ClearFormatting()
$r2.Find.Execute("enhance", $true, $false, $false, $false, $false, $true, 1, $false, "broaden", 2)

# --- Edit 3: paragraph "........accuracy" -> "enhance.......accuracy" (note: one fewer ellipsis char)
$p3 = $d.Paragraphs.Item(12)
$r3 = $p3.Range
$r3.Find.ClearFormatting()
$r3.Find.Execute("………………………..accuracy", $true, $false, $false, $false, $false, $true, 1, $false, "enhance……………………..accuracy", 2)

# --- Edit 4: paragraph "...strengthen........rusty / unused skills" -> "... brush off........rusty / unused skills"
$p4 = $d.Paragraphs.Item(13)
$r4 = $p4.Range
$r4.Find.ClearFormatting()
$r4.Find.Execute("strengthen", $true, $false, $false, $false, $false, $true, 1, $false, " brush off", 2)
